# Update the "Greeting" value for rule R10 (cell E8) and leave that cell
# selected/active, matching the author's "update file with jgit" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
